$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 50
$ws.Range("B4").Value = 40
$ws.Range("E6").Select()
